# Apply updated sector/portfolio metrics following a refreshed data pull.
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet: Portfolio Raw
# ------------------------------------------------------------------
$wsRaw = $wb.Worksheets.Item("Portfolio Raw")

$wsRaw.Range("F2").Value = 19.266619
$wsRaw.Range("G2").Value = 0.023426486394750368

$wsRaw.Range("C3").Value = 0.0050999997
$wsRaw.Range("F3").Value = 24.338146
$wsRaw.Range("G3").Value = 0.0007899811224648264

$wsRaw.Range("F4").Value = 25.276854
$wsRaw.Range("G4").Value = 0.005807474659536062

$wsRaw.Range("C5").Value = 0.0073
$wsRaw.Range("F5").Value = 30.535118
$wsRaw.Range("G5").Value = 0.008136970526567382

# ------------------------------------------------------------------
# Sheet: Sector Reference
# ------------------------------------------------------------------
$wsSector = $wb.Worksheets.Item("Sector Reference")

$wsSector.Range("F2").Value = 19.266619
$wsSector.Range("G2").Value = 0.023426486394750368

$wsSector.Range("C3").Value = 0.0049
$wsSector.Range("F3").Value = 19.435196
$wsSector.Range("G3").Value = 0.023494077359232035

$wsSector.Range("C4").Value = 0.0404
$wsSector.Range("F4").Value = 12.234955
$wsSector.Range("G4").Value = 0.010959942259112204

$wsSector.Range("F5").Value = 11.503075
$wsSector.Range("G5").Value = 0.015276398741521398

$wsSector.Range("C6").Value = 0.0376
$wsSector.Range("F6").Value = 7.943082
$wsSector.Range("G6").Value = 0.012469486233323755

$wsSector.Range("F7").Value = 114.134834
$wsSector.Range("G7").Value = 0.016341055995351877

$wsSector.Range("C8").Value = 0.005
$wsSector.Range("F8").Value = 17.287224
$wsSector.Range("G8").Value = 0.026088105144970536

$wsSector.Range("C9").Value = 0.0098
$wsSector.Range("F9").Value = 14.699733
$wsSector.Range("G9").Value = 0.009215757137694345

$wsSector.Range("C10").Value = 0.0106
$wsSector.Range("F10").Value = 21.459576
$wsSector.Range("G10").Value = 0.0111764705882353

$wsSector.Range("C11").Value = 0.0544
$wsSector.Range("F11").Value = 8.801818
$wsSector.Range("G11").Value = 0.004319210201563161

$wsSector.Range("F12").Value = 59.92975
$wsSector.Range("G12").Value = 0.000895563516120112

$wsSector.Range("F13").Value = 25.276854
$wsSector.Range("G13").Value = 0.005807474659536062

$wsSector.Range("F14").Value = 50.269135
$wsSector.Range("G14").Value = 0.0027861524743400964

$wsSector.Range("C15").Value = 0.0072000003
$wsSector.Range("F15").Value = 26.17453
$wsSector.Range("G15").Value = 0.022030313711667252

$wsSector.Range("C16").Value = 0.0062
$wsSector.Range("F16").Value = 39.091465
$wsSector.Range("G16").Value = 0.006816421378776107

$wsSector.Range("C17").Value = 0.0381
$wsSector.Range("F17").Value = 8.503472
$wsSector.Range("G17").Value = 0.0032276051820767683

$wsSector.Range("C18").Value = 0.0168
$wsSector.Range("F18").Value = 7.285276
$wsSector.Range("G18").Value = 0.00461022632020115

$wsSector.Range("C19").Value = 0.0146
$wsSector.Range("F19").Value = 22.547329
$wsSector.Range("G19").Value = 0.006114186565314343

$wsSector.Range("F20").Value = 28.958546
$wsSector.Range("G20").Value = 0.005976357267951003

$wsSector.Range("C21").Value = 0.0095
$wsSector.Range("F21").Value = 48.527775
$wsSector.Range("G21").Value = 0.016451895408482572

$wsSector.Range("C22").Value = 0.0619
$wsSector.Range("F22").Value = 9.1712475
$wsSector.Range("G22").Value = 0.010492700729927026

$wsSector.Range("C23").Value = 0.0091
$wsSector.Range("F23").Value = 21.673786
$wsSector.Range("G23").Value = 0.0010739216037228417

$wsSector.Range("F24").Value = -69.26923
$wsSector.Range("G24").Value = 0.000554938956714651

$wsSector.Range("C25").Value = 0.01453913045652174
$wsSector.Range("F25").Value = 22.822002065217397
$wsSector.Range("G25").Value = 0.010417597731764565

$wsSector.Range("C26").Value = 0.01803216302069066
$wsSector.Range("F26").Value = 30.46340248938333
$wsSector.Range("G26").Value = 0.007663708478974685

# ------------------------------------------------------------------
# Sheet: Portfolio Normalized
# ------------------------------------------------------------------
$wsNorm = $wb.Worksheets.Item("Portfolio Normalized")

$wsNorm.Range("B2").Value = 96.05319707489454

$wsNorm.Range("C3").Value = 50
$wsNorm.Range("D3").Value = 50

$wsNorm.Range("B4").Value = 99.22009806722488
$wsNorm.Range("C4").Value = 59.497378817906416
$wsNorm.Range("D4").Value = 71.82444586443386

$wsNorm.Range("C5").Value = 60.200539477765204
$wsNorm.Range("D5").Value = 72.56788595322104

# ------------------------------------------------------------------
# Sheet: Z-Score Comparison
# ------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("Z-Score Comparison")

$wsZ.Range("C2").Value = -0.796818046550493
$wsZ.Range("K2").Value = 19.266619
$wsZ.Range("L2").Value = -0.1167099790135547
$wsZ.Range("M2").Formula = "=STANDARDIZE(19.266619,AVERAGE('Sector Reference'!F:F),STDEV('Sector Reference'!F:F))"
$wsZ.Range("N2").Value = 0.023426486394750368
$wsZ.Range("O2").Value = 1.6974665331641425
$wsZ.Range("P2").Formula = "=STANDARDIZE(0.023426486394750368,AVERAGE('Sector Reference'!G:G),STDEV('Sector Reference'!G:G))"

$wsZ.Range("B3").Value = 0.0050999997
$wsZ.Range("C3").Value = -0.7913866647836603
$wsZ.Range("D3").Formula = "=STANDARDIZE(0.0050999997,AVERAGE('Sector Reference'!C:C),STDEV('Sector Reference'!C:C))"
$wsZ.Range("K3").Value = 24.338146
$wsZ.Range("L3").Value = 0.04976935637150137
$wsZ.Range("M3").Formula = "=STANDARDIZE(24.338146,AVERAGE('Sector Reference'!F:F),STDEV('Sector Reference'!F:F))"
$wsZ.Range("N3").Value = 0.0007899811224648264
$wsZ.Range("O3").Value = -1.2562608084210165
$wsZ.Range("P3").Formula = "=STANDARDIZE(0.0007899811224648264,AVERAGE('Sector Reference'!G:G),STDEV('Sector Reference'!G:G))"

$wsZ.Range("C4").Value = -0.8891518027251564
$wsZ.Range("K4").Value = 25.276854
$wsZ.Range("L4").Value = 0.08058364247520063
$wsZ.Range("M4").Formula = "=STANDARDIZE(25.276854,AVERAGE('Sector Reference'!F:F),STDEV('Sector Reference'!F:F))"
$wsZ.Range("N4").Value = 0.005807474659536062
$wsZ.Range("O4").Value = -0.6015525101034744
$wsZ.Range("P4").Formula = "=STANDARDIZE(0.005807474659536062,AVERAGE('Sector Reference'!G:G),STDEV('Sector Reference'!G:G))"

$wsZ.Range("B5").Value = 0.0073
$wsZ.Range("C5").Value = -0.6718958911468715
$wsZ.Range("D5").Formula = "=STANDARDIZE(0.0073,AVERAGE('Sector Reference'!C:C),STDEV('Sector Reference'!C:C))"
$wsZ.Range("K5").Value = 30.535118
$wsZ.Range("L5").Value = 0.2531928578060402
$wsZ.Range("M5").Formula = "=STANDARDIZE(30.535118,AVERAGE('Sector Reference'!F:F),STDEV('Sector Reference'!F:F))"
$wsZ.Range("N5").Value = 0.008136970526567382
$wsZ.Range("O5").Value = -0.29758793819650925
$wsZ.Range("P5").Formula = "=STANDARDIZE(0.008136970526567382,AVERAGE('Sector Reference'!G:G),STDEV('Sector Reference'!G:G))"
